$d = $word.ActiveDocument

# The three inline logo pictures (two "PearsonLogo" copies in the footers and
# one "BTec_Logo-Orange" copy in the header) were exported with their
# wp:docPr / pic:cNvPr "name" attributes swapped relative to their actual
# media part. Fix the displayed picture names:
#   - PearsonLogo pictures: image2.png -> image1.png
#   - BTec_Logo-Orange picture: image1.jpg -> image2.jpg
# These "name" values aren't exposed as a settable InlineShape property in
# the Word object model (InlineShape has no .Name), so edit the underlying
# WordprocessingML directly via Document.WordOpenXML, which round-trips the
# full OOXML package as a single string.

$xml = $d.WordOpenXML

$xml = $xml.Replace('name="image2.png"', 'name="image1.png"')
$xml = $xml.Replace('name="image1.jpg"', 'name="image2.jpg"')

$d.WordOpenXML = $xml
